$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "Datos actualizados" timestamp (A1): 21:05 -> 21:35
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 21:35"

# ---------------------------------------------------------------------------
# 2) Refresh the figures for "Estados Unidos" (row 4)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 1420421
$ws.Range("C4").Value = 11785
$ws.Range("D4").Value = 305723
$ws.Range("E4").Value = 1030345
$ws.Range("F4").Value = 16373
$ws.Range("G4").Value = 928
$ws.Range("H4").Value = 84353

# ---------------------------------------------------------------------------
# 3) "Sierra Leona" moves up (new data) and pushes "Estado de Palestina" /
#    "Republica del Chad" one row down (rows 130-132). "Congo" (row 133)
#    is unaffected.
# ---------------------------------------------------------------------------
$ws.Range("A130").Value = "Sierra Leona"
$ws.Range("B130").Value = 387
$ws.Range("C130").Value = 49
$ws.Range("D130").Value = 97
$ws.Range("E130").Value = 264
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 7
$ws.Range("H130").Value = 26

$ws.Range("A131").Value = "Estado de Palestina"
$ws.Range("B131").Value = 375
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 310
$ws.Range("E131").Value = 63
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 2

$ws.Range("A132").Value = "Republica del Chad"
$ws.Range("B132").Value = 372
$ws.Range("C132").Value = 15
$ws.Range("D132").Value = 78
$ws.Range("E132").Value = 252
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 42

# Row 133 ("Congo") stays exactly as-is (333, 0, 53, 269, 0, 0, 11).

# ---------------------------------------------------------------------------
# 4) "Nueva Caledonia" moves up ahead of "Belice" (rows 193-194), the
#    underlying per-country figures travel with their country name.
# ---------------------------------------------------------------------------
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 18
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 16
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2
